# Update Metadata sheet timestamp (A2)
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 10:17 AM"

# Update "Industry Analysis" sheet: refresh the "1 Year" column (F) values for rows 2-76
$wsIndustry = $wb.Worksheets.Item("Industry Analysis")
$industryUpdates = @(
  @{Row=2; Val=21.0016},
  @{Row=3; Val=-16.2396},
  @{Row=4; Val=27.1317},
  @{Row=5; Val=-50.6494},
  @{Row=6; Val=53.2813},
  @{Row=7; Val=-8.106199999999999},
  @{Row=8; Val=-9.552099999999999},
  @{Row=9; Val=36.3756},
  @{Row=10; Val=-6.1314},
  @{Row=11; Val=31.9081},
  @{Row=12; Val=-18.4955},
  @{Row=13; Val=14.0155},
  @{Row=14; Val=-36.0718},
  @{Row=15; Val=-0.1622},
  @{Row=16; Val=0.1459},
  @{Row=17; Val=-22.0012},
  @{Row=18; Val=1.0561},
  @{Row=19; Val=-27.708},
  @{Row=20; Val=47.7309},
  @{Row=21; Val=12.0959},
  @{Row=22; Val=95.1491},
  @{Row=23; Val=-50.2657},
  @{Row=24; Val=-13.3427},
  @{Row=25; Val=-9.9316},
  @{Row=26; Val=5.8244},
  @{Row=27; Val=-32.7692},
  @{Row=28; Val=-24.8224},
  @{Row=29; Val=-18.4191},
  @{Row=30; Val=25.8569},
  @{Row=31; Val=58.4712},
  @{Row=32; Val=-3.3862},
  @{Row=33; Val=-6.3282},
  @{Row=34; Val=27.7203},
  @{Row=35; Val=4.4873},
  @{Row=36; Val=-4.9458},
  @{Row=37; Val=3.6074},
  @{Row=38; Val=-23.3973},
  @{Row=39; Val=8.7355},
  @{Row=40; Val=-5.8541},
  @{Row=41; Val=-8.3934},
  @{Row=42; Val=20.3818},
  @{Row=43; Val=14.3164},
  @{Row=44; Val=-12.6846},
  @{Row=45; Val=28.4075},
  @{Row=46; Val=-1.1135},
  @{Row=47; Val=-37.1997},
  @{Row=48; Val=-29.8569},
  @{Row=49; Val=-27.5511},
  @{Row=50; Val=-49.7478},
  @{Row=51; Val=-51.8002},
  @{Row=52; Val=-38.5254},
  @{Row=53; Val=-12.4886},
  @{Row=54; Val=-5.0725},
  @{Row=55; Val=-17.7445},
  @{Row=56; Val=-26.636},
  @{Row=57; Val=-29.3361},
  @{Row=58; Val=-11.9574},
  @{Row=59; Val=-24.5687},
  @{Row=60; Val=-12.3},
  @{Row=61; Val=-10.9446},
  @{Row=62; Val=-17.1229},
  @{Row=63; Val=-9.5038},
  @{Row=64; Val=54.2749},
  @{Row=65; Val=-43.4736},
  @{Row=66; Val=13.2687},
  @{Row=67; Val=12.7149},
  @{Row=68; Val=24.8057},
  @{Row=69; Val=-17.0328},
  @{Row=70; Val=-6.8927},
  @{Row=71; Val=13.6034},
  @{Row=72; Val=3.9995},
  @{Row=73; Val=-16.226},
  @{Row=74; Val=-16.2448},
  @{Row=75; Val=28.6924},
  @{Row=76; Val=48.9752}
)
foreach ($u in $industryUpdates) {
    $wsIndustry.Cells.Item($u.Row, 6).Value = $u.Val
}

# Update "Stock List" sheet: the top row (CAPTRU-RE1) dropped out of the list,
# every other row shifts up by one, and a new row (TRAVELFOOD) is appended at the bottom.
$wsStock = $wb.Worksheets.Item("Stock List")
$wsStock.Rows.Item(2).Delete()

$wsStock.Range("A76").Value = "📋"
$wsStock.Range("B76").Value = "TRAVELFOOD"
$wsStock.Range("C76").Value = "TRAVELFOOD"
$wsStock.Range("D76").Value = 1316.3
$wsStock.Range("E76").Value = 0.1141
$wsStock.Range("F76").Value = "N/A"
$wsStock.Range("G76").Value = "N/A"
$wsStock.Range("H76").Value = 17332.9705
